$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2054380664652568
$ws.Range("C2").Value = 0.552870090634441
$ws.Range("J2").Value = 0.00906344410876133
$ws.Range("O2").Value = 0.003021148036253776
$ws.Range("P2").Value = 0.1389728096676737
$ws.Range("S2").Value = 0.09063444108761329
$ws.Range("B3").Value = 0.01047120418848168
$ws.Range("C3").Value = 0.03664921465968586
$ws.Range("J3").Value = 0.03664921465968586
$ws.Range("P3").Value = 0.680628272251309
$ws.Range("S3").Value = 0.2356020942408377
$ws.Range("B6").Value = 0.06808510638297872
$ws.Range("D6").Value = 0.01276595744680851
$ws.Range("F6").Value = 0.03404255319148936
$ws.Range("J6").Value = 0.2425531914893617
$ws.Range("O6").Value = 0.02127659574468085
$ws.Range("Q6").Value = 0.1148936170212766
$ws.Range("R6").Value = 0.08936170212765958
$ws.Range("S6").Value = 0.4170212765957447
$ws.Range("B7").Value = 0.1105990783410138
$ws.Range("D7").Value = 0.004608294930875576
$ws.Range("F7").Value = 0.05529953917050692
$ws.Range("J7").Value = 0.09216589861751152
$ws.Range("O7").Value = 0.0184331797235023
$ws.Range("Q7").Value = 0.1751152073732719
$ws.Range("R7").Value = 0.09216589861751152
$ws.Range("S7").Value = 0.4516129032258064
$ws.Range("B8").Value = 0.08834586466165413
$ws.Range("D8").Value = 0.009398496240601503
$ws.Range("F8").Value = 0.07894736842105263
$ws.Range("J8").Value = 0.112781954887218
$ws.Range("O8").Value = 0.005639097744360902
$ws.Range("Q8").Value = 0.2048872180451128
$ws.Range("R8").Value = 0.06578947368421052
$ws.Range("S8").Value = 0.4342105263157895
$ws.Range("B9").Value = 0.12
$ws.Range("D9").Value = 0.006666666666666667
$ws.Range("F9").Value = 0.05333333333333334
$ws.Range("J9").Value = 0.08
$ws.Range("O9").Value = 0.01333333333333333
$ws.Range("Q9").Value = 0.16
$ws.Range("R9").Value = 0.08
$ws.Range("S9").Value = 0.4866666666666667
$ws.Range("B10").Value = 0.115919629057187
$ws.Range("D10").Value = 0.0115919629057187
$ws.Range("F10").Value = 0.06800618238021638
$ws.Range("J10").Value = 0.1027820710973725
$ws.Range("O10").Value = 0.01468315301391035
$ws.Range("Q10").Value = 0.2102009273570325
$ws.Range("R10").Value = 0.07727975270479134
$ws.Range("S10").Value = 0.3995363214837713
$ws.Range("G11").Value = 0.1218130311614731
$ws.Range("J11").Value = 0.1104815864022663
$ws.Range("K11").Value = 0.1784702549575071
$ws.Range("L11").Value = 0.5864022662889519
$ws.Range("S11").Value = 0.0028328611898017
$ws.Range("G12").Value = 0.7109004739336493
$ws.Range("J12").Value = 0.2654028436018958
$ws.Range("L12").Value = 0.004739336492890996
$ws.Range("S12").Value = 0.01895734597156398
$ws.Range("G13").Value = 0.6530612244897959
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.06122448979591837
$ws.Range("F15").Value = 0.01376146788990826
$ws.Range("H15").Value = 0.1880733944954129
$ws.Range("I15").Value = 0.05504587155963303
$ws.Range("J15").Value = 0.3532110091743119
$ws.Range("K15").Value = 0.05963302752293578
$ws.Range("M15").Value = 0.02293577981651376
$ws.Range("O15").Value = 0.06880733944954129
$ws.Range("S15").Value = 0.2385321100917431
$ws.Range("F16").Value = 0.02162162162162162
$ws.Range("H16").Value = 0.2216216216216216
$ws.Range("I16").Value = 0.05405405405405406
$ws.Range("J16").Value = 0.4108108108108108
$ws.Range("K16").Value = 0.1027027027027027
$ws.Range("M16").Value = 0.01621621621621622
$ws.Range("O16").Value = 0.03243243243243243
$ws.Range("S16").Value = 0.1405405405405405
$ws.Range("F17").Value = 0.02155172413793104
$ws.Range("H17").Value = 0.1831896551724138
$ws.Range("I17").Value = 0.0625
$ws.Range("J17").Value = 0.418103448275862
$ws.Range("K17").Value = 0.125
$ws.Range("M17").Value = 0.01724137931034483
$ws.Range("N17").Value = 0.002155172413793103
$ws.Range("O17").Value = 0.04956896551724138
$ws.Range("S17").Value = 0.1206896551724138
$ws.Range("F18").Value = 0.02150537634408602
$ws.Range("H18").Value = 0.2043010752688172
$ws.Range("I18").Value = 0.05376344086021505
$ws.Range("J18").Value = 0.3602150537634409
$ws.Range("K18").Value = 0.08602150537634409
$ws.Range("M18").Value = 0.02688172043010753
$ws.Range("O18").Value = 0.08064516129032258
$ws.Range("S18").Value = 0.1666666666666667
$ws.Range("F19").Value = 0.01798063623789765
$ws.Range("H19").Value = 0.2282157676348548
$ws.Range("I19").Value = 0.06154910096818811
$ws.Range("J19").Value = 0.3457814661134163
$ws.Range("K19").Value = 0.1230982019363762
$ws.Range("M19").Value = 0.02143845089903181
$ws.Range("N19").Value = 0.001383125864453665
$ws.Range("O19").Value = 0.06777316735822959
$ws.Range("S19").Value = 0.1327800829875519
